# Insert a new row above the current row 2 on the "Schedules" sheet,
# shifting the existing schedule rows down by one, and populate the
# new row's first cell with the "default" label (new shared string).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedules")

$ws.Rows.Item(2).Insert() | Out-Null
$ws.Range("A2").Value = "default"

# Update the sheet's active cell / selection to match the author's
# final cursor position.
$ws.Range("H13").Select() | Out-Null
